$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8302608728408813
$ws.Range("B1").Value = 1.245790362358093
$ws.Range("C1").Value = 2.250858068466187
$ws.Range("D1").Value = 2.353326559066772
$ws.Range("E1").Value = 1.956781506538391
